# Applies scheduled market-data refresh values to the Malboro_Profits workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 43
$ws.Range("H43").Value = 5256.375
$ws.Range("J43").Value = 5507.4287
$ws.Range("L43").Value = 5507.4287
$ws.Range("N43").Value = -5645.4287
# Row 51
$ws.Range("H51").Value = 13000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
# Row 62
$ws.Range("H62").Value = 3996.6667
$ws.Range("I62").Value = 3990
$ws.Range("K62").Value = 3990
$ws.Range("M62").Value = -3366
# Row 65
$ws.Range("H65").Value = 3996.6667
$ws.Range("I65").Value = 3990
$ws.Range("K65").Value = 19950
$ws.Range("M65").Value = -16830
# Row 69
$ws.Range("H69").Value = 13069
$ws.Range("J69").Value = 13069
$ws.Range("L69").Value = 39207
$ws.Range("N69").Value = -40955
# Row 72
$ws.Range("H72").Value = 13069
$ws.Range("J72").Value = 13069
$ws.Range("L72").Value = 117621
$ws.Range("N72").Value = -126357
# Row 76
$ws.Range("H76").Value = 24309.467
$ws.Range("J76").Value = 24212.25
$ws.Range("L76").Value = 24212.25
$ws.Range("N76").Value = -24842.25
# Row 79
$ws.Range("H79").Value = 24309.467
$ws.Range("J79").Value = 24212.25
$ws.Range("L79").Value = 24212.25
$ws.Range("N79").Value = -26396.25
# Row 111
$ws.Range("H111").Value = 1763.8334
$ws.Range("I111").Value = 1264.5
$ws.Range("J111").Value = 2013.5
$ws.Range("K111").Value = 3793.5
$ws.Range("L111").Value = 6040.5
$ws.Range("M111").Value = -726.5
$ws.Range("N111").Value = -12174.5
# Row 132
$ws.Range("H132").Value = 19798.045
$ws.Range("I132").Value = 12252.556
$ws.Range("J132").Value = 53752.75
$ws.Range("K132").Value = 36757.66800000001
$ws.Range("L132").Value = 161258.25
$ws.Range("M132").Value = -34227.66800000001
$ws.Range("N132").Value = -166318.25
# Row 138
$ws.Range("H138").Value = 2112.9285
$ws.Range("J138").Value = 3174.762
$ws.Range("L138").Value = 9524.286
$ws.Range("N138").Value = -19804.286

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 4810697
$ws.Range("I32").Value = 6580752
$ws.Range("J32").Value = 6262.2144
$ws.Range("K32").Value = 6580752
$ws.Range("L32").Value = 6262.2144
$ws.Range("M32").Value = -6580465
$ws.Range("N32").Value = -6836.2144
# Row 61
$ws.Range("H61").Value = 909044.75
$ws.Range("I61").Value = 3169.72
$ws.Range("K61").Value = 3169.72
$ws.Range("M61").Value = -2957.72
# Row 74
$ws.Range("H74").Value = 14910.409
$ws.Range("I74").Value = 1426.7667
$ws.Range("J74").Value = 43803.93
$ws.Range("K74").Value = 1426.7667
$ws.Range("L74").Value = 43803.93
$ws.Range("M74").Value = -552.7666999999999
$ws.Range("N74").Value = -45551.93
# Row 77
$ws.Range("H77").Value = 14910.409
$ws.Range("I77").Value = 1426.7667
$ws.Range("J77").Value = 43803.93
$ws.Range("K77").Value = 7133.8335
$ws.Range("L77").Value = 219019.65
$ws.Range("M77").Value = -2765.8335
$ws.Range("N77").Value = -227755.65
# Row 136
$ws.Range("H136").Value = 909044.75
$ws.Range("I136").Value = 3169.72
$ws.Range("K136").Value = 9509.16
$ws.Range("M136").Value = -6959.16

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 5
$ws.Range("H5").Value = 314.5
$ws.Range("I5").Value = 314.5
$ws.Range("K5").Value = 314.5
$ws.Range("M5").Value = -201.5
# Row 12
$ws.Range("H12").Value = 926.6667
$ws.Range("I12").Value = 120
$ws.Range("K12").Value = 120
$ws.Range("M12").Value = 48
# Row 19
$ws.Range("H19").Value = 105932
$ws.Range("J19").Value = 105932
$ws.Range("L19").Value = 105932
$ws.Range("N19").Value = -106278
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("L46").ClearContents()
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("L76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("L79").ClearContents()
# Row 87
$ws.Range("H87").Value = 250000
$ws.Range("J87").Value = 250000
$ws.Range("L87").Value = 250000
$ws.Range("N87").Value = -252496
# Row 90
$ws.Range("H90").Value = 250000
$ws.Range("J90").Value = 250000
$ws.Range("L90").Value = 750000
$ws.Range("N90").Value = -762480

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 22
$ws.Range("H22").Value = 351.8889
$ws.Range("I22").Value = 351.125
$ws.Range("J22").Value = 358
$ws.Range("K22").Value = 351.125
$ws.Range("L22").Value = 358
$ws.Range("M22").Value = -1.125
$ws.Range("N22").Value = -1058
# Row 23
$ws.Range("H23").Value = 4500
$ws.Range("I23").Value = 4500
$ws.Range("K23").Value = 4500
$ws.Range("M23").Value = -4260
# Row 27
$ws.Range("H27").Value = 4500
$ws.Range("I27").Value = 4500
$ws.Range("K27").Value = 4500
$ws.Range("M27").Value = -4308
# Row 41
$ws.Range("H41").Value = 1754.5
$ws.Range("I41").Value = 1754.5
$ws.Range("K41").Value = 1754.5
$ws.Range("M41").Value = -1326.5
# Row 132
$ws.Range("H132").Value = 31197038
$ws.Range("I132").Value = 2537.56
$ws.Range("K132").Value = 7612.68
$ws.Range("M132").Value = -5082.68
# Row 141
$ws.Range("H141").Value = 173418
$ws.Range("I141").Value = 35000
$ws.Range("J141").Value = 196487.67
$ws.Range("K141").Value = 35000
$ws.Range("L141").Value = 196487.67
$ws.Range("M141").Value = -29820
$ws.Range("N141").Value = -206847.67

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 131
$ws.Range("H131").Value = 1477.33
$ws.Range("J131").Value = 1479.0928
$ws.Range("L131").Value = 4437.278399999999
$ws.Range("N131").Value = -14517.2784
# Row 137
$ws.Range("H137").Value = 9566
$ws.Range("I137").Value = 5999.5
$ws.Range("K137").Value = 17998.5
$ws.Range("M137").Value = -12898.5
# Row 138
$ws.Range("H138").Value = 5399.1
$ws.Range("I138").Value = 3665.1667
$ws.Range("K138").Value = 10995.5001
$ws.Range("M138").Value = -5855.500100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 68
$ws.Range("H68").Value = 50147
$ws.Range("J68").Value = 50295
$ws.Range("L68").Value = 50295
$ws.Range("N68").Value = -51917
# Row 71
$ws.Range("H71").Value = 50147
$ws.Range("J71").Value = 50295
$ws.Range("L71").Value = 150885
$ws.Range("N71").Value = -158997
# Row 80
$ws.Range("H80").Value = 1000
$ws.Range("J80").Value = 1000
$ws.Range("L80").Value = 1000
$ws.Range("N80").Value = -2996
# Row 83
$ws.Range("H83").Value = 1000
$ws.Range("J83").Value = 1000
$ws.Range("L83").Value = 5000
$ws.Range("N83").Value = -14984
# Row 88
$ws.Range("H88").Value = 95000
$ws.Range("I88").Value = 95000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 95000
$ws.Range("N88").Value = 0
$ws.Range("M88").Value = -94549
$ws.Range("L88").ClearContents()
# Row 91
$ws.Range("H91").Value = 95000
$ws.Range("I91").Value = 95000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 95000
$ws.Range("N91").Value = 0
$ws.Range("M91").Value = -93440
$ws.Range("L91").ClearContents()
# Row 123
$ws.Range("H123").Value = 55199.125
$ws.Range("J123").Value = 55199.125
$ws.Range("L123").Value = 55199.125
$ws.Range("N123").Value = -60099.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 22
$ws.Range("H22").Value = 4998.3335
# Row 27
$ws.Range("H27").Value = 4998.3335
# Row 46
$ws.Range("H46").Value = 850
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 100
$ws.Range("H100").Value = 3436.6667
$ws.Range("J100").Value = 4397.5
$ws.Range("L100").Value = 4397.5
$ws.Range("N100").Value = -5479.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 6
$ws.Range("H6").Value = 2601.6667
$ws.Range("I6").Value = 105
$ws.Range("J6").Value = 3101
$ws.Range("K6").Value = 105
$ws.Range("L6").Value = 3101
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = -3331
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("L48").ClearContents()
# Row 96
$ws.Range("H96").Value = 1930.6154
$ws.Range("I96").Value = 2081.3333
$ws.Range("J96").Value = 1885.4
$ws.Range("K96").Value = 2081.3333
$ws.Range("L96").Value = 1885.4
$ws.Range("M96").Value = -708.3332999999998
$ws.Range("N96").Value = -4631.4
# Row 136
$ws.Range("H136").Value = 401832.38
$ws.Range("I136").Value = 1955.3334
$ws.Range("K136").Value = 5866.0002
$ws.Range("M136").Value = -3316.0002
